$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'After Searched'
$ws.Range("B1").Value = 'After Pool Filter'
$ws.Range("C1").Value = 'After Top Reviewed'
$ws.Range("A2").Value = 'YEZA INN
Opened in 2025
8.5/10
Very good2 reviews
Near Yaya CentreShow on Map
2-bed Room
x4
Free Cancellation
Breakfast included
Only 1 left at this price
Special Discount
10% off
₹ 6,163
₹ 5,489
Total price: ₹ 32,440
1 room × 5 nights incl. taxes & fees
Check Availability
Sign in for member prices'
$ws.Range("B2").Value = 'Woodmere Serviced Apartment
8.6/10
Very good68 reviews
Near Yaya CentreShow on Map
Standard Two-Bedroom Apartment
x4
Entire unit 59㎡
2 bedrooms
2 beds
₹ 5,239
Total price: ₹ 30,881
1 room × 5 nights incl. taxes & fees
Check Availability'
$ws.Range("C2").Value = 'Holiday Inn NAIROBI TWO RIVERS MALL by IHG
9.7/10
Amazing94 reviews
"Great location"
"Great service"
Near Village MarketShow on Map
No. 13 of 4-Star Select Hotels in Nairobi
Family Room
x4
Free Cancellation
₹ 19,986
Total price: ₹ 127,909
1 room × 5 nights incl. taxes & fees
Check Availability'
$ws.Range("A3").Value = 'Javilla Eagles Safari Guest house
7.7/10
Good34 reviews
8.4 km from centreShow on Map
Comfort Apartment, 2 Bedrooms, Non Smoking, Ground Floor
x4
Entire unit 55㎡
2 bedrooms
2 beds
Only 1 left at this price
₹ 5,058
Total price: ₹ 29,337
1 room × 5 nights incl. taxes & fees
Check Availability'
$ws.Range("B3").Value = 'Kenya Comfort Suites
7.9/10
Good11 reviews
Near Yaya CentreShow on Map
Standard Quadruple Room
x4
₹ 4,785
Total price: ₹ 29,391
1 room × 5 nights incl. taxes & fees
Check Availability'
$ws.Range("C3").Value = 'JW Marriott Hotel Nairobi
9.8/10
Outstanding46 reviews
"Great stay!"
"Great service"
Near National Museum of KenyaShow on Map
No. 1 of Luxury Hotels in Nairobi
3 Bedroom Apartment, Bedroom 1: 1 King, Bedroom 2: 1 King, Bedroom 3: 2 Doubles
x4
Free Cancellation
Breakfast included
₹ 203,399
Total price: ₹ 1,284,055
1 room × 5 nights incl. taxes & fees
Check Availability'
$ws.Range("A4").Value = 'Woodmere Serviced Apartment
8.6/10
Very good68 reviews
Near Yaya CentreShow on Map
Standard Two-Bedroom Apartment
x4
Entire unit 59㎡
2 bedrooms
2 beds
₹ 5,239
Total price: ₹ 30,881
1 room × 5 nights incl. taxes & fees
Check Availability'
$ws.Range("B4").Value = 'The King Post
7.7/10
Good17 reviews
Near The Sarit Expo CentreShow on Map
Three-Bedroom Apartment
x6
Entire unit
3 bedrooms
4 beds
₹ 4,459
Total price: ₹ 26,263
1 room × 5 nights incl. taxes & fees
Check Availability'
$ws.Range("C4").Value = 'Yaya Hotel & Apartments
9.4/10
Amazing49 reviews
"Clean and tidy"
"Great service"
Near Yaya CentreShow on Map
No. 19 of 4-Star Select Hotels in Nairobi
Two-Bedroom Luxury Apartment
x4
Entire unit 125㎡
2 bedrooms
2 beds
Free Cancellation
Earn ₹ 962.26 in Trip Coins
Last booked 14 hrs ago
Limited Time Offer
15% off
₹ 19,632
₹ 16,304
Total price: ₹ 96,194
1 room × 5 nights incl. taxes & fees
Check Availability
Sign in for member prices'
$ws.Range("A5").Value = 'Kenya Comfort Suites
7.9/10
Good11 reviews
Near Yaya CentreShow on Map
Standard Quadruple Room
x4
₹ 4,785
Total price: ₹ 29,391
1 room × 5 nights incl. taxes & fees
Check Availability'
$ws.Range("B5").Value = 'Eldon Apartments & Suites
7.1/10
5 reviews
Near Wilson AirportShow on Map
Standard Two-Bedroom Apartment
x4
Entire apartment 18㎡
2 bedrooms
3 beds
Special Discount
₹ 112 off
₹ 4,441
₹ 4,329
Total price: ₹ 25,553
1 room × 5 nights incl. taxes & fees
Check Availability'
$ws.Range("C5").Value = 'Pan Pacific Serviced Suites Nairobi
9.4/10
Amazing48 reviews
"Great rooms"
"Great location"
Near National Museum of KenyaShow on Map
No. 13 of Premium Hotels in Nairobi
Two Bedroom Suite King & Twin
x4
Entire unit 108㎡
2 bedrooms
2 beds
Last booked 3 hrs ago
₹ 23,235
Total price: ₹ 145,219
1 room × 5 nights incl. taxes & fees
Check Availability'
$ws.Range("A6").Value = 'Kester International Apartment Hotel
Opened in 2025
9.8/10
Outstanding28 reviews
"Great service"
"Great rooms"
Near Yaya CentreShow on Map
Boutique 2-bedroom And 1-living Room Suite
x4
Entire apartment 95㎡
2 bedrooms
3 beds
Earn ₹ 547.03 in Trip Coins
Only 3 left at this price
Special Discount
11% off
₹ 7,141
₹ 6,284
Total price: ₹ 36,444
1 room × 5 nights incl. taxes & fees
Check Availability
Sign in for member prices'
$ws.Range("B6").Value = 'Lux Suites Riara One Residency Angama
Opened in 2025
9.9/10
Outstanding17 reviews
"Clean and tidy"
"Great location"
Near Yaya CentreShow on Map
Family Room
x4
Entire apartment 98㎡
2 bedrooms
2 beds
Free Cancellation
Breakfast included
Earn ₹ 1,323.34 in Trip Coins
Only 5 left at this price
Special Discount
8% off
₹ 16,398
₹ 14,946
Total price: ₹ 88,178
1 room × 5 nights incl. taxes & fees
Check Availability
Sign in for member prices'
$ws.Range("C6").Value = 'Palacina the Residence & the Suites
9.2/10
Great48 reviews
Near Yaya CentreShow on Map
No. 7 of Premium Hotels in Nairobi
2 Bedroom Executive Penthouse
x4
Entire unit 168㎡
2 bedrooms
3 beds
Free Cancellation
Breakfast included
Earn ₹ 2,506.77 in Trip Coins
Only 1 left at this price
Special Discount
20% off
₹ 36,674
₹ 28,320
Total price: ₹ 167,085
1 room × 5 nights incl. taxes & fees
Check Availability
Sign in for member prices'
$ws.Range("A7").Value = 'Holiday Inn NAIROBI TWO RIVERS MALL by IHG
Ad
9.7/10
Amazing94 reviews
"Great location"
"Great service"
Near Village MarketShow on Map
No. 13 of 4-Star Select Hotels in Nairobi
Family Room
x4
Free Cancellation
₹ 19,986
Total price: ₹ 127,909
1 room × 5 nights incl. taxes & fees
Check Availability'
$ws.Range("B7").Value = 'Holiday Inn NAIROBI TWO RIVERS MALL by IHG
Ad
9.7/10
Amazing94 reviews
"Great location"
"Great service"
Near Village MarketShow on Map
No. 13 of 4-Star Select Hotels in Nairobi
Family Room
x4
Free Cancellation
₹ 19,986
Total price: ₹ 127,909
1 room × 5 nights incl. taxes & fees
Check Availability'
$ws.Range("C7").Value = 'Executive Residency by Best Western Nairobi
9.2/10
Great46 reviews
Near The Sarit Expo CentreShow on Map
Two Bedroom Apartment
x4
Breakfast included
₹ 18,186
Total price: ₹ 111,845
1 room × 5 nights incl. taxes & fees
Check Availability'
$ws.Range("A8").Value = 'Eldon Apartments & Suites
7.1/10
5 reviews
Near Wilson AirportShow on Map
Standard Two-Bedroom Apartment
x4
Entire apartment 18㎡
2 bedrooms
3 beds
Special Discount
₹ 112 off
₹ 4,441
₹ 4,329
Total price: ₹ 25,553
1 room × 5 nights incl. taxes & fees
Check Availability'
$ws.Range("B8").Value = 'Arcadia Hotel
Renovated in 2025
8.8/10
Very good13 reviews
Near Yaya CentreShow on Map
Two-Bedroom Suite
x4
Entire unit 110㎡
2 bedrooms
2 beds
Free Cancellation
Only 5 left at this price
₹ 7,704
Total price: ₹ 44,682
1 room × 5 nights incl. taxes & fees
Check Availability
Sign in for member prices'
$ws.Range("C8").Value = 'Mövenpick Hotel & Residences Nairobi
9.0/10
Great97 reviews
"Great service"
"Great location"
Near The Sarit Expo CentreShow on Map
No. 1 of Gourmet Hotels in Nairobi
Two-Bedroom Residence
x4
Free Cancellation
₹ 18,849
Total price: ₹ 117,805
1 room × 5 nights incl. taxes & fees
Check Availability'
$ws.Range("A9").Value = 'The King Post
7.7/10
Good17 reviews
Near The Sarit Expo CentreShow on Map
Three-Bedroom Apartment
x6
Entire unit
3 bedrooms
4 beds
₹ 4,459
Total price: ₹ 26,263
1 room × 5 nights incl. taxes & fees
Check Availability'
$ws.Range("B9").Value = 'Lavington Residences By Trianum
Opened in 2025
9.6/10
Amazing7 reviews
Near Yaya CentreShow on Map
Executive Two-Bedroom Apartment
x4
Only 1 left at this price
₹ 9,562
Total price: ₹ 56,417
1 room × 5 nights incl. taxes & fees
Check Availability
Sign in for member prices'
$ws.Range("C9").Value = 'Fairview Hotel Nairobi, Vignette Collection by IHG
Renovated in 2025
9.0/10
Great49 reviews
Near Giraffe manorShow on Map
No. 17 of 4-Star Select Hotels in Nairobi
Fairview Suite
x4
Free Cancellation
₹ 54,648
Total price: ₹ 349,747
1 room × 5 nights incl. taxes & fees
Check Availability'
$ws.Range("A10").Value = 'Maskan Suites
8.9/10
Very good34 reviews
Near Yaya CentreShow on Map
Superior Apartment, 2 Bedrooms, Private Bathroom, City View
x4
Entire apartment
2 bedrooms
2 beds
Member deal
₹ 6,945
Total price: ₹ 40,941
1 room × 5 nights incl. taxes & fees
Check Availability'
$ws.Range("B10").Value = 'Yaya Hotel & Apartments
9.4/10
Amazing49 reviews
"Clean and tidy"
"Great service"
Near Yaya CentreShow on Map
No. 19 of 4-Star Select Hotels in Nairobi
Two-Bedroom Luxury Apartment
x4
Entire unit 125㎡
2 bedrooms
2 beds
Free Cancellation
Earn ₹ 962.26 in Trip Coins
Last booked 14 hrs ago
Limited Time Offer
15% off
₹ 19,632
₹ 16,304
Total price: ₹ 96,194
1 room × 5 nights incl. taxes & fees
Check Availability
Sign in for member prices'
$ws.Range("C10").Value = 'Windsor Golf Hotel & Country Club
New to Trip.com
8.8/10
Very good64 reviews
9.7 km from centreShow on Map
No. 6 of 4-Star Select Hotels in Nairobi
Two- Bedroom Cottage
x4
Entire unit 65㎡
1 bedroom
4 beds
Breakfast included
₹ 36,895
Total price: ₹ 217,678
1 room × 5 nights incl. taxes & fees
Check Availability'
$ws.Range("A11").Value = 'Lux Suites Riara One Residency Angama
Opened in 2025
9.9/10
Outstanding17 reviews
"Clean and tidy"
"Great location"
Near Yaya CentreShow on Map
Family Room
x4
Entire apartment 98㎡
2 bedrooms
2 beds
Free Cancellation
Breakfast included
Earn ₹ 1,323.34 in Trip Coins
Only 5 left at this price
Special Discount
8% off
₹ 16,398
₹ 14,946
Total price: ₹ 88,178
1 room × 5 nights incl. taxes & fees
Check Availability
Sign in for member prices'
$ws.Range("B11").Value = 'Mövenpick Hotel & Residences Nairobi
9.0/10
Great97 reviews
"Great service"
"Great location"
Near The Sarit Expo CentreShow on Map
No. 1 of Gourmet Hotels in Nairobi
Two-Bedroom Residence
x4
Free Cancellation
₹ 18,849
Total price: ₹ 117,805
1 room × 5 nights incl. taxes & fees
Check Availability'
$ws.Range("C11").Value = 'Four Points by Sheraton Nairobi Hurlingham
8.8/10
Very good60 reviews
"Great service"
"Delicious breakfast"
Near Giraffe manorShow on Map
Executive Suite
x4
₹ 29,260
Total price: ₹ 187,305
1 room × 5 nights incl. taxes & fees
Check Availability'
$ws.Range("A12").Value = 'Pullman Nairobi Upper Hill
Renovated in 2025
Ad
9.4/10
Amazing60 reviews
Near Giraffe manorShow on Map
No. 3 of 4-Star Select Hotels in Nairobi
2
Deluxe Room With Two Double Beds
Earn ₹ 1,697.05 in Trip Coins
Last booked 2 hrs ago
₹ 13,576
Total price: ₹ 169,703
2 rooms × 5 nights incl. taxes & fees
Check Availability'
$ws.Range("B12").Value = 'Pullman Nairobi Upper Hill
Renovated in 2025
Ad
9.4/10
Amazing60 reviews
Near Giraffe manorShow on Map
No. 3 of 4-Star Select Hotels in Nairobi
2
Deluxe Room With Two Double Beds
Earn ₹ 1,697.05 in Trip Coins
Last booked 2 hrs ago
₹ 13,576
Total price: ₹ 169,703
2 rooms × 5 nights incl. taxes & fees
Check Availability'
$ws.Range("A13").Value = 'Arcadia Hotel
Renovated in 2025
8.8/10
Very good13 reviews
Near Yaya CentreShow on Map
Two-Bedroom Suite
x4
Entire unit 110㎡
2 bedrooms
2 beds
Free Cancellation
Only 5 left at this price
₹ 7,704
Total price: ₹ 44,682
1 room × 5 nights incl. taxes & fees
Check Availability
Sign in for member prices'
